$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AG2").Value = "dichotic_phase"
$ws.Range("AJ2").Value = 7
$ws.Range("AG3").Value = "MAB_and_AFACT"
$ws.Range("AJ3").Value = 3
$ws.Range("AG4").Value = "Dichotic_and_AFACT"
$ws.Range("AJ4").Value = 2
$ws.Range("AG5").Value = "Dichotic_and_AFACT"
$ws.Range("AJ5").Value = 3
$ws.Range("AG6").Value = "MAB_phase"
$ws.Range("AJ6").Value = 6
$ws.Range("AG7").Value = "MAB_phase"
$ws.Range("AJ8").Value = 1
$ws.Range("AG9").Value = "dichotic_phase"
$ws.Range("AJ9").Value = 4
$ws.Range("AG10").Value = "dichotic_phase"
$ws.Range("AJ10").Value = 5
$ws.Range("AG11").Value = "MAB_and_Digit_after"
$ws.Range("AJ11").Value = 8
$ws.Range("AG12").Value = "Dichotic_and_AFACT"
$ws.Range("AJ12").Value = 4
$ws.Range("AG13").Value = "Dichotic_and_AFACT"
$ws.Range("AJ13").Value = 2
$ws.Range("AG14").Value = "MAB_and_AFACT"
$ws.Range("AJ14").Value = 1
$ws.Range("AJ15").Value = 7
$ws.Range("AG16").Value = "MAB_and_Digit_after"
$ws.Range("AJ16").Value = 8
$ws.Range("AG17").Value = "MAB_and_AFACT"
$ws.Range("AJ17").Value = 5
$ws.Range("AG18").Value = "Dichotic_and_AFACT"
$ws.Range("AG19").Value = "Dichotic_and_AFACT"
$ws.Range("AJ19").Value = 7
$ws.Range("AG20").Value = "dichotic_phase"
$ws.Range("AJ20").Value = 1
$ws.Range("AG21").Value = "Digit_before_and_AFACT"
$ws.Range("AJ21").Value = 3
$ws.Range("AG22").Value = "MAB_phase"
$ws.Range("AJ22").Value = 7
$ws.Range("AG23").Value = "MAB_phase"
$ws.Range("AJ23").Value = 1
$ws.Range("AG24").Value = "Dichotic_and_AFACT"
$ws.Range("AJ24").Value = 4
$ws.Range("AJ25").Value = 6
$ws.Range("AG26").Value = "MAB_and_Digit_after"
$ws.Range("AJ26").Value = 4
$ws.Range("AG27").Value = "Digit_before_and_AFACT"
$ws.Range("AJ27").Value = 8
$ws.Range("AG28").Value = "dichotic_phase"
$ws.Range("AJ28").Value = 4
$ws.Range("AG29").Value = "MAB_and_AFACT"
$ws.Range("AJ29").Value = 1
$ws.Range("AG30").Value = "MAB_and_AFACT"
$ws.Range("AJ30").Value = 5
$ws.Range("AG31").Value = "MAB_and_Digit_after"
$ws.Range("AJ31").Value = 7
$ws.Range("AG32").Value = "dichotic_phase"
$ws.Range("AJ32").Value = 1
$ws.Range("AJ33").Value = 3
$ws.Range("AJ34").Value = 3
$ws.Range("AG35").Value = "Digit_before_and_AFACT"
$ws.Range("AJ35").Value = 7
$ws.Range("AG36").Value = "dichotic_phase"
$ws.Range("AJ36").Value = 5
$ws.Range("AG37").Value = "dichotic_phase"
$ws.Range("AG38").Value = "Dichotic_and_AFACT"
$ws.Range("AJ38").Value = 7
$ws.Range("AJ39").Value = 4
$ws.Range("AG40").Value = "dichotic_phase"
$ws.Range("AJ40").Value = 4
$ws.Range("AG41").Value = "Dichotic_and_AFACT"
$ws.Range("AJ41").Value = 2
$ws.Range("AG42").Value = "Dichotic_and_AFACT"
$ws.Range("AJ42").Value = 6
$ws.Range("AG43").Value = "dichotic_phase"
$ws.Range("AJ43").Value = 2
$ws.Range("AG44").Value = "MAB_and_Digit_after"
$ws.Range("AJ44").Value = 6
$ws.Range("AG45").Value = "MAB_and_Digit_after"
$ws.Range("AJ45").Value = 8
$ws.Range("AG46").Value = "dichotic_phase"
$ws.Range("AJ46").Value = 2
$ws.Range("AG47").Value = "MAB_and_Digit_after"
$ws.Range("AJ47").Value = 5
$ws.Range("AG48").Value = "dichotic_phase"
$ws.Range("AJ48").Value = 5
$ws.Range("AG49").Value = "Dichotic_and_AFACT"
$ws.Range("AJ49").Value = 3
$ws.Range("AJ50").Value = 5
$ws.Range("AG51").Value = "Digit_before_and_AFACT"
$ws.Range("AJ51").Value = 6
$ws.Range("AJ52").Value = 7
$ws.Range("AG53").Value = "dichotic_phase"
$ws.Range("AJ53").Value = 8
$ws.Range("AG55").Value = "MAB_phase"
$ws.Range("AJ55").Value = 6
$ws.Range("AG56").Value = "dichotic_phase"
$ws.Range("AJ56").Value = 8
$ws.Range("AG57").Value = "Dichotic_and_AFACT"
$ws.Range("AJ57").Value = 1
$ws.Range("AG58").Value = "Dichotic_and_AFACT"
$ws.Range("AJ58").Value = 5
$ws.Range("AG59").Value = "Digit_before_and_AFACT"
$ws.Range("AJ59").Value = 1
$ws.Range("AG60").Value = "MAB_and_AFACT"
$ws.Range("AJ60").Value = 7
$ws.Range("AG61").Value = "MAB_phase"
$ws.Range("AJ61").Value = 3
$ws.Range("AG62").Value = "MAB_and_Digit_after"
$ws.Range("AJ62").Value = 2
$ws.Range("AG63").Value = "Dichotic_and_AFACT"
$ws.Range("AJ63").Value = 7
$ws.Range("AG64").Value = "MAB_phase"
$ws.Range("AJ64").Value = 6
$ws.Range("AG65").Value = "MAB_phase"
$ws.Range("AJ65").Value = 3
$ws.Range("AG66").Value = "Digit_before_and_AFACT"
$ws.Range("AG67").Value = "MAB_and_AFACT"
$ws.Range("AJ67").Value = 3
$ws.Range("AG68").Value = "Dichotic_and_AFACT"
$ws.Range("AJ68").Value = 5
$ws.Range("AG69").Value = "MAB_and_Digit_after"
$ws.Range("AJ69").Value = 7
$ws.Range("AG70").Value = "Digit_before_and_AFACT"
$ws.Range("AJ70").Value = 4
$ws.Range("AG71").Value = "dichotic_phase"
$ws.Range("AJ71").Value = 6
$ws.Range("AG72").Value = "MAB_and_AFACT"
$ws.Range("AJ72").Value = 7
$ws.Range("AG73").Value = "MAB_and_AFACT"
$ws.Range("AJ73").Value = 4
$ws.Range("AG74").Value = "dichotic_phase"
$ws.Range("AJ74").Value = 8
$ws.Range("AG75").Value = "Dichotic_and_AFACT"
$ws.Range("AJ75").Value = 8
$ws.Range("AG76").Value = "dichotic_phase"
$ws.Range("AJ76").Value = 1
$ws.Range("AJ77").Value = 2
$ws.Range("AJ78").Value = 5
$ws.Range("AG79").Value = "MAB_and_AFACT"
$ws.Range("AJ79").Value = 7
$ws.Range("AG80").Value = "Digit_before_and_AFACT"
$ws.Range("AJ80").Value = 6
$ws.Range("AG81").Value = "Dichotic_and_AFACT"
